# Updated MCH102 to MCH251 -- add the MCH241 archive record as row 2 of
# the collections sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the new record row (row 2) -----------------------------
# Columns: A=identifier  C=title  D=(blank) E=levelOfDescription
#          F=extentAndMedium  G=notes  H=(blank)
# (B, and I..Z intentionally stay empty, like the header's unused columns.)
$ws.Range("A2").Value = "MCH241"
$ws.Range("C2").Value = "CORRESPONDENCE BETWEEN ANDRE ODENDAAL & MR SIMON EGERT MR. SIMON TO ANDRE ODENDAAL, MY WEEDS WERE FLOURISHING"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"
$ws.Range("H2").Value = ""

# --- Match the body-row font used for the rest of the sheet ----------
# (Calibri 10pt, automatic/theme text color) on every cell of the new
# row that carries data (or is an intentional blank placeholder).
foreach ($addr in @("A2", "C2", "D2", "E2", "F2", "G2", "H2")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

# --- Re-establish the frozen header row / selection on row 2 ---------
$ws.Range("A2:H2").Select()
$excel.ActiveWindow.FreezePanes = $true
